{"js": "// The document ends with a paragraph \"Xnj-nj\" that carries the hidden\n// \"_GoBack\" bookmark (an empty/collapsed range right after the text, before\n// the paragraph mark). We need to append a new paragraph \"Df,h\" after it,\n// and leave the \"_GoBack\" bookmark collapsed at the very end of the\n// document (i.e. right after \"Df,h\"), exactly like Word does when new text\n// is typed at the end of a document.\n\nconst body = context.document.body;\n\n// Drop the existing \"_GoBack\" bookmark first so that inserting a new one\n// with the same name later doesn't just add a second, stale bookmark.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Find the last paragraph and append the new paragraph right after it.\nbody.paragraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = body.paragraphs.items[body.paragraphs.items.length - 1];\nlastParagraph.insertParagraph(\"Df,h\", Word.InsertLocation.after);\nawait context.sync();\n\n// Re-create the \"_GoBack\" bookmark collapsed at the new end of the body,\n// matching where Word leaves it after the most recent edit.\nconst endRange = body.getRange(\"End\");\nendRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "# The document ends with a paragraph \"Xnj-nj\" that carries the hidden\n# \"_GoBack\" bookmark, collapsed right after the text (before the paragraph\n# mark). We need to append a new paragraph \"Df,h\" after it, and leave the\n# \"_GoBack\" bookmark collapsed at the very end of the document (i.e. right\n# after \"Df,h\"), exactly like Word does when new text is typed at the end\n# of a document.\n\n$d = $word.ActiveDocument\n\n# Add a new, empty paragraph after the current last paragraph.\n$lastParagraph = $d.Paragraphs.Last\n$lastParagraph.Range.InsertParagraphAfter()\n\n# The freshly inserted paragraph is now the last one; its Range.Start is\n# where the new text belongs.\n$newParagraph = $d.Paragraphs.Last\n$insertPos = $newParagraph.Range.Start\n\n# Type the new paragraph's text plus a one-character throw-away marker.\n# (Planting a *collapsed* bookmark exactly at a paragraph's closing\n# boundary, i.e. right before its paragraph mark, lands it in the wrong\n# spot on this host, so we temporarily keep a harmless character after the\n# insertion point, anchor the bookmark just before that character -\n# a safe, non-boundary position - and delete the marker afterwards. The\n# bookmark stays put once anchored.)\n$newText = \"Df,h\"\n$marker = \"~\"\n$insertRange = $d.Range($insertPos, $insertPos)\n$insertRange.InsertAfter($newText + $marker)\n\n# Re-anchor \"_GoBack\" collapsed right after the real text, before the marker.\n$bookmarkPos = $insertPos + $newText.Length\n$bmRange = $d.Range($bookmarkPos, $bookmarkPos)\n$d.Bookmarks.Add(\"_GoBack\", $bmRange)\n\n# Remove the throw-away marker now that the bookmark is safely anchored.\n$markerRange = $d.Range($bookmarkPos, $bookmarkPos + $marker.Length)\n$markerRange.Delete()\n"}
